# Apply the grading-sheet update: adjust a handful of raw score entries,
# change the per-row total formula's weight for column I from 3.5 to 4,
# mark a few rows as "переписаны верно все номера", and add a small
# scratch area below the table (labels for an "i"/"j"/"arr" snippet).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Raw input corrections (columns B..I) ---------------------------------
$ws.Range("B3").Value = 3
$ws.Range("D3").Value = 3
$ws.Range("E3").Value = -2

$ws.Range("G5").Value = 5

$ws.Range("E7").Value = 4

$ws.Range("E8").Value = -2

$ws.Range("E16").Value = 4
$ws.Range("G16").Value = 0

$ws.Range("F20").Value = 4
$ws.Range("G20").Value = 5

$ws.Range("E22").Value = 4
$ws.Range("F22").Value = 4
$ws.Range("G22").Value = 5

$ws.Range("F23").Value = 4

# --- Total-score formula: weight for column I changes 3.5 -> 4 -----------
$ws.Range("J2:J25").Formula = "=MAX(B2+2*C2+D2+E2+G2+H2+2*F2+4*I2,0)"

# --- Notes: mark these rows as fully/correctly rewritten -----------------
$ws.Range("N20").Value = "переписаны верно все номера"
$ws.Range("N22").Value = "переписаны верно все номера"
$ws.Range("N23").Value = "переписаны верно все номера"

# --- New scratch area below the table -------------------------------------
$ws.Range("A32").Value = "i"
$ws.Range("A33").Value = "j"
$ws.Range("A30").Value = "arr"

$ws.Range("B30:I38").Select()
